# LMDI: Versjon 1.0.6 6e05e801b0fc67a31d9121f33125496b6f7ed95a
#
# - "Experimental" row (B7 on the Metadata sheet) loses its "false" value
#   (cell becomes blank, and the now-unused "false" shared string drops out
#   of the shared-string table, shifting every later index down by one).
# - "Date" row (B8) value changes from 2025-04-11 to 2025-09-12. The new
#   value must stay a literal text string (not get auto-parsed into a date
#   serial number) and keep the cell's original style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Clear the "Experimental" value cell entirely -> removes the "false"
# shared string and collapses later indices, matching the diff.
$ws.Cells.Item(7, 2).Value = ""

# Write the new date as literal text, without ever touching NumberFormat
# on a kept cell (that would otherwise leave a permanent unused style
# behind). Build the text in a scratch cell via a quoted formula (so Excel
# never runs its "looks like a date" autodetection on it), copy just the
# computed VALUE onto B8 (so B8 keeps its existing style/s="2"), then wipe
# the scratch cell.
$scratch = $ws.Cells.Item(20, 5)
$scratch.Formula = "=""2025-09-12"""
$scratch.Copy()
$ws.Cells.Item(8, 2).PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = $false
